# Results for alpha = 0.01 -- fill in column E (the "alpha = 0.01" series)
# on Sheet1 with the joint model's total-training-cost values, then move
# the active selection to M21 (matching the author's final cursor spot).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(
    26565.7975693,
    24377.3680291,
    24560.7594295,
    24453.998558,
    24744.0533018,
    24549.5441608,
    24553.8706493,
    23936.5731478,
    24190.8561316,
    22814.6516104,
    23255.7587996,
    23231.3569994,
    23160.6402683,
    24986.9566936,
    22978.0805016,
    23113.1382852,
    23313.9190412,
    22632.9259491,
    23201.073472,
    22742.8499775
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 5).Value = $values[$i]
}

$ws.Range("M21").Select()
